$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.598.36'
$ws.Range('D3').Value = '2.460.15'
$ws.Range('E3').Value = '  -0.43%  '
$ws.Range('D5').Value = '318.37'
$ws.Range('E5').Value = '  +0.67%  '
$ws.Range('D6').Value = '91.26'
$ws.Range('E6').Value = '  -1.57%  '
$ws.Range('D7').Value = '0.547'
$ws.Range('E7').Value = '  -1.11%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('E9').Value = '  -2.37%  '
$ws.Range('E10').Value = '  -3.09%  '
$ws.Range('D11').Value = '32.61'
$ws.Range('E11').Value = '  -0.75%  '
$ws.Range('E12').Value = '  -0.84%  '
$ws.Range('D13').Value = '2.838.52'
$ws.Range('E13').Value = '  -0.50%  '
$ws.Range('D14').Value = '6.83'
$ws.Range('E14').Value = '  -1.14%  '
$ws.Range('E15').Value = '  -2.25%  '
$ws.Range('D16').Value = '2.447.92'
$ws.Range('E16').Value = '  -0.98%  '
$ws.Range('D17').Value = '0.784'
$ws.Range('E17').Value = '  -0.13%  '
$ws.Range('D18').Value = '41.491.02'
$ws.Range('E18').Value = '  -0.38%  '
$ws.Range('D19').Value = '6.38'
$ws.Range('E19').Value = '  -1.56%  '
$ws.Range('D20').Value = '0.0₃0936'
$ws.Range('E20').Value = '  -3.61%  '
$ws.Range('D21').Value = '71.74'
$ws.Range('E21').Value = '  +0.60%  '
$ws.Range('D22').Value = '11.08'
$ws.Range('E22').Value = '  -3.26%  '
$ws.Range('D23').Value = '238.08'
$ws.Range('E23').Value = '  -0.42%  '
$ws.Range('E24').Value = '  +0.13%  '
$ws.Range('D25').Value = '1.92'
$ws.Range('E25').Value = '  +0.37%  '
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('D27').Value = '24.58'
$ws.Range('E27').Value = '  -0.54%  '
$ws.Range('E28').Value = '  -1.52%  '
$ws.Range('D29').Value = '9.65'
$ws.Range('E29').Value = '  -1.82%  '
$ws.Range('D30').Value = '36.03'
$ws.Range('E30').Value = '  +1.62%  '
$ws.Range('D31').Value = '157.41'
$ws.Range('E31').Value = '  +0.88%  '
$ws.Range('D32').Value = '5.39'
$ws.Range('E32').Value = '  -2.32%  '
$ws.Range('E33').Value = '  +0.07%  '
$ws.Range('E34').Value = '  -0.42%  '
$ws.Range('D35').Value = '0.0759'
$ws.Range('E35').Value = '  -0.89%  '
$ws.Range('D36').Value = '16.86'
$ws.Range('E36').Value = '  -4.06%  '
$ws.Range('B37').Value = 'Stellar'
$ws.Range('C37').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D37').Value = '0.116'
$ws.Range('E37').Value = '  +0.27%  '
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').Value = '2.88'
$ws.Range('E38').Value = '  -0.29%  '
$ws.Range('D39').Value = '1.82'
$ws.Range('E39').Value = '  +0.58%  '
$ws.Range('E40').Value = '  -0.03%  '
$ws.Range('D41').Value = '3.99'
$ws.Range('E41').Value = '  -0.20%  '
$ws.Range('E42').Value = '  -7.24%  '
$ws.Range('D43').Value = '1.991.04'
$ws.Range('E43').Value = '  +1.13%  '
$ws.Range('E44').Value = '  -1.65%  '
$ws.Range('D45').Value = '18.43'
$ws.Range('E45').Value = '  -2.23%  '
$ws.Range('E46').Value = '  -0.66%  '
$ws.Range('D47').Value = '9.47'
$ws.Range('E47').Value = '  +4.41%  '
$ws.Range('D48').Value = '2.718.11'
$ws.Range('E48').Value = '  +0.43%  '
$ws.Range('B49').Value = 'BitcoinSV'
$ws.Range('C49').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D49').Value = '75.55'
$ws.Range('E49').Value = '  +4.18%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').Value = '96.83'
$ws.Range('E50').Value = '  -0.49%  '
$ws.Range('D51').Value = '66.52'
$ws.Range('E51').Value = '  -0.55%  '
